$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old C2 value (Player 2 score row becomes empty / removed)
$ws.Range("C2").ClearContents()

# Update row 2 values: A2 = 0, B2 = 26
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 26

# Add row 3: A3 = 28, B3 = 0
$ws.Range("A3").Value = 28
$ws.Range("B3").Value = 0

# Add row 4: A4 = 0, B4 = 16
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 16
